# Add a new "BodyText" paragraph after the last existing paragraph
# ("And another one from the deploy keys repo"), containing three
# separate runs:
#   1. "SSH_DEPLOY_KEY updated (now without a new line at the end)"
#   2. " "
#   3. "Now using the correct deploy key"

$d = $word.ActiveDocument

$last = $d.Paragraphs.Last

# Step 1: create three temporary paragraphs (same style as the
# preceding paragraph, i.e. "BodyText"), one run of text each.
# Using InsertParagraphAfter()+InsertAfter() on distinct paragraphs
# keeps each chunk of text in its own run (runs only get coalesced
# together when inserted back-to-back into the very same paragraph).
$last.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.InsertAfter("SSH_DEPLOY_KEY updated (now without a new line at the end)")

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.InsertAfter(" ")

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.InsertAfter("Now using the correct deploy key")

# Step 2: merge the three temporary paragraphs into a single paragraph
# by deleting the paragraph marks between them. This keeps the three
# pieces of text as three distinct runs inside one paragraph (exactly
# matching the target diff) instead of merging their text together.
$countAfterInsert = $d.Paragraphs.Count
$firstIdx = $countAfterInsert - 2

$firstPara = $d.Paragraphs.Item($firstIdx)
$mark1 = $firstPara.Range.End - 1
$d.Range($mark1, $mark1 + 1).Delete()

$firstPara = $d.Paragraphs.Item($firstIdx)
$mark2 = $firstPara.Range.End - 1
$d.Range($mark2, $mark2 + 1).Delete()

$final = $d.Paragraphs.Item($firstIdx)
Write-Output $final.Range.Text
